# Change fonts in the "short report" Word template.
# The custom "CorpoS" font (and, for the List Bullet style, the default
# "Arial" font) used by the report's paragraph styles is replaced with
# "Calibri" everywhere it appears.

$d = $word.ActiveDocument

$styleNames = @(
    "Metadata",
    "Risk Bold List",
    "Risk Italics List",
    "Measure List",
    "Risk Italics",
    "Measure Indent",
    "List Bullet"
)

foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    $style.Font.Name = "Calibri"
}
